# Update cryptos list data (prices and volume percentages) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

$ws.Range("D2").Value = "56.740.02"
$ws.Range("E2").Value = "  +10.49%  "

$ws.Range("D3").Value = "3.253.95"
$ws.Range("E3").Value = "  +5.77%  "

$ws.Range("E4").Value = "  +0.08%  "

Set-TextValue "D5" "395.58"
$ws.Range("E5").Value = "  +0.31%  "

Set-TextValue "D6" "109.26"
$ws.Range("E6").Value = "  +6.62%  "

Set-TextValue "D7" "0.559"
$ws.Range("E7").Value = "  +4.58%  "

$ws.Range("E8").Value = "  +0.00%  "

Set-TextValue "D9" "0.623"
$ws.Range("E9").Value = "  +6.14%  "

Set-TextValue "D10" "39.16"
$ws.Range("E10").Value = "  +4.51%  "

Set-TextValue "D11" "0.0956"
$ws.Range("E11").Value = "  +12.12%  "

$ws.Range("E12").Value = "  +2.16%  "

$ws.Range("D13").Value = "3.778.23"
$ws.Range("E13").Value = "  +6.11%  "

Set-TextValue "D14" "8.19"
$ws.Range("E14").Value = "  +6.37%  "

Set-TextValue "D15" "19.07"
$ws.Range("E15").Value = "  +2.05%  "

$ws.Range("D16").Value = "3.270.24"
$ws.Range("E16").Value = "  +6.63%  "

$ws.Range("E17").Value = "  +0.98%  "

Set-TextValue "D18" "10.79"
$ws.Range("E18").Value = "  +2.22%  "

$ws.Range("D19").Value = "56.670.25"
$ws.Range("E19").Value = "  +10.34%  "

Set-TextValue "D20" "3.29"
$ws.Range("E20").Value = "  +3.99%  "

$ws.Range("E21").Value = "  +9.26%  "

Set-TextValue "D22" "12.88"
$ws.Range("E22").Value = "  +4.39%  "

Set-TextValue "D23" "304.29"
$ws.Range("E23").Value = "  +14.97%  "

Set-TextValue "D24" "74.95"
$ws.Range("E24").Value = "  +6.68%  "

Set-TextValue "D25" "3.14"
$ws.Range("E25").Value = "  -2.11%  "

Set-TextValue "D26" "28.07"
$ws.Range("E26").Value = "  +3.97%  "

$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D27" "4.39"
$ws.Range("E27").Value = "  +4.96%  "

$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D28" "7.85"
$ws.Range("E28").Value = "  -0.31%  "

$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D29" "0.170"
$ws.Range("E29").Value = "  +3.23%  "

$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D30" "7.24"
$ws.Range("E30").Value = "  +0.63%  "

$ws.Range("E31").Value = "  -0.40%  "

Set-TextValue "D32" "0.109"
$ws.Range("E32").Value = "  +4.61%  "

$ws.Range("E33").Value = "  +2.73%  "

Set-TextValue "D34" "37.31"
$ws.Range("E34").Value = "  +2.41%  "

Set-TextValue "D35" "0.0479"
$ws.Range("E35").Value = "  -2.77%  "

Set-TextValue "D36" "2.13"
$ws.Range("E36").Value = "  +3.15%  "

Set-TextValue "D37" "51.45"
$ws.Range("E37").Value = "  +3.07%  "

Set-TextValue "D38" "3.52"
$ws.Range("E38").Value = "  +5.68%  "

$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D39" "0.999"
$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D40" "3.08"
$ws.Range("E40").Value = "  +22.00%  "

Set-TextValue "D41" "135.20"
$ws.Range("E41").Value = "  +4.86%  "

Set-TextValue "D42" "1.92"
$ws.Range("E42").Value = "  +4.17%  "

Set-TextValue "D43" "17.22"
$ws.Range("E43").Value = "  +3.31%  "

$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D44" "4.00"
$ws.Range("E44").Value = "  -0.68%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D45" "0.120"
$ws.Range("E45").Value = "  +4.09%  "

$ws.Range("E46").Value = "  -2.99%  "

Set-TextValue "D47" "22.02"
$ws.Range("E47").Value = "  +1.25%  "

$ws.Range("D48").Value = "2.144.19"
$ws.Range("E48").Value = "  +3.48%  "

$ws.Range("E49").Value = "  +2.23%  "

$ws.Range("E50").Value = "  -5.93%  "

Set-TextValue "D51" "2.00"
$ws.Range("E51").Value = "  +37.64%  "
